# OLX Monitor run @ 2026-02-15 21:16 — append a fresh snapshot row to each
# profile-detail sheet and refresh the summary sheet's timestamps.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2026-02-15 21:16"

# ---------------------------------------------------------------------
# 1) PODSUMOWANIE (summary) sheet: bump the "last checked" timestamps and
#    zero out the "Nowe (+)" counters now that this run has recorded them.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("PODSUMOWANIE")
$summary.Range("B2").Value = $newTimestamp
$summary.Range("B3").Value = $newTimestamp
$summary.Range("B4").Value = $newTimestamp
$summary.Range("B5").Value = $newTimestamp
$summary.Range("B6").Value = $newTimestamp
$summary.Range("D4").Value = 0
$summary.Range("D5").Value = 0

# ---------------------------------------------------------------------
# 2) Detail sheets: append row 3 with this run's results.
# ---------------------------------------------------------------------
function Add-MonitorRow {
    param($ws, $total, $newDetails, $removedDetails, $newIds)

    # Row height matches the other data rows.
    $ws.Rows.Item(3).RowHeight = 18

    # Pull in the per-column formatting used for the sheet's other rows
    # before the values are written, so we inherit the correct style ids.
    $ws.Range("A2").Copy()
    $ws.Range("A3").PasteSpecial(-4122)
    $ws.Range("H2").Copy()
    $ws.Range("H3").PasteSpecial(-4122)

    # B/C/D/E use the plain (non-highlighted) numeric style — pull that
    # from the PODSUMOWANIE sheet, which already carries it.
    $summary.Range("B2").Copy()
    $ws.Range("B3").PasteSpecial(-4122)
    $ws.Range("C3").PasteSpecial(-4122)
    $ws.Range("D3").PasteSpecial(-4122)
    $summary.Range("D2").Copy()
    $ws.Range("E3").PasteSpecial(-4122)

    # F/G ("Szczegóły nowych/usuniętych") use a left-aligned, unshaded
    # variant of the existing style — seed from the plain numeric style
    # and flip the alignment to left.
    $summary.Range("B2").Copy()
    $ws.Range("F3").PasteSpecial(-4122)
    $ws.Range("F3").HorizontalAlignment = -4131
    $ws.Range("F3").Copy()
    $ws.Range("G3").PasteSpecial(-4122)

    $ws.Range("A3").Value = $newTimestamp
    $ws.Range("B3").Value = $total
    $ws.Range("C3").Value = 0
    $ws.Range("D3").Value = 0
    $ws.Range("E3").Value = 0
    $ws.Range("F3").Value = $newDetails
    $ws.Range("G3").Value = $removedDetails
    $ws.Range("H3").Value = "OK"
    if ($newIds -ne "") {
        $ws.Range("I3").Value = $newIds
    }

    $ws.Application.CutCopyMode = $false
}

$wsLublin = $wb.Worksheets.Item("wszystkie-lublin")
Add-MonitorRow $wsLublin 431 "—" "—" ""

$wsArtymiuk = $wb.Worksheets.Item("artymiuk")
Add-MonitorRow $wsArtymiuk 0 "—" "—" ""

$wsPoqui = $wb.Worksheets.Item("poqui")
Add-MonitorRow $wsPoqui 5 "—" "—" "18KAEc|183ger|17NeTz|1951OR|17vbYq"

$wsStylowe = $wb.Worksheets.Item("stylowepokoje")
Add-MonitorRow $wsStylowe 2 "—" "—" "16ZeYm|195dLc"

$wsVillahome = $wb.Worksheets.Item("villahome")
Add-MonitorRow $wsVillahome 0 "—" "—" ""
